$d = $word.ActiveDocument

$replacements = @(
    @("306÷4=76, 2", "985÷9=109, 4"),
    @("743÷6=123, 5", "663÷7=94, 5"),
    @("360÷4=90, 0", "835÷8=104, 3"),
    @("923÷2=461, 1", "424÷2=212, 0"),
    @("473÷2=236, 1", "333÷4=83, 1"),
    @("984÷4=246, 0", "647÷8=80, 7"),
    @("278÷9=30, 8", "493÷2=246, 1"),
    @("230÷3=76, 2", "771÷6=128, 3"),
    @("429÷9=47, 6", "548÷8=68, 4"),
    @("604÷2=302, 0", "781÷3=260, 1"),
    @("335÷2=167, 1", "769÷2=384, 1"),
    @("728÷3=242, 2", "870÷5=174, 0"),
    @("197÷7=28, 1", "419÷5=83, 4"),
    @("562÷9=62, 4", "387÷2=193, 1"),
    @("124÷4=31, 0", "213÷7=30, 3"),
    @("151÷4=37, 3", "807÷5=161, 2"),
    @("301÷9=33, 4", "529÷4=132, 1"),
    @("576÷3=192, 0", "779÷6=129, 5"),
    @("817÷5=163, 2", "335÷7=47, 6"),
    @("857÷6=142, 5", "925÷8=115, 5"),
    @("948÷2=474, 0", "447÷8=55, 7"),
    @("676÷4=169, 0", "584÷7=83, 3"),
    @("556÷4=139, 0", "843÷9=93, 6"),
    @("458÷6=76, 2", "730÷6=121, 4"),
    @("593÷3=197, 2", "207÷2=103, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done replacing $($replacements.Count) entries"
